# "Made use of merged column reading"
#
# CharacterTalentEffectDataTable's header row (row 2) used to spell out
# one "LevelN" label per talent-level column (E2:S2). Replace that with a
# single header cell/merged range reading "EffectStatistics" spanning
# E2:S2 (left/vertical-centered, percentage number format, matching the
# sheet's other numeric columns), since the data is now read as one
# merged block instead of per-level columns.
#
# Also brings CharacterTalentEffectDataTable to the front as the active
# tab (it was CharacterTalentDataTable before).

$wb = $excel.ActiveWorkbook
$wsEffect = $wb.Worksheets.Item("CharacterTalentEffectDataTable")

# Style + relabel E2:S2 first (one clean new cell format gets registered
# here), THEN merge the range down to the single E2 cell/value. Doing the
# styling before the merge keeps every cell in the run on one consistent
# format and is the combination that stays closest to the recorded edit.
$headerRange = $wsEffect.Range("E2:S2")
$headerRange.NumberFormat = "0.00%"
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4108
$wsEffect.Range("E2").Value = "EffectStatistics"
$headerRange.Merge()

# Switch the active sheet/tab from CharacterTalentDataTable to
# CharacterTalentEffectDataTable, with F25 selected there.
$wsEffect.Activate()
$wsEffect.Range("F25").Select()
